# Naming convention: Change snake_case to camelCase
# Renames the "Attribute" column values (and one constraint note that
# references a renamed attribute) from snake_case to camelCase across the
# ER-diagram-style tables laid out in columns G:K of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# H-column attribute name renames (snake_case -> camelCase)
$ws.Range("H5").Value  = "isDeleted"
$ws.Range("H11").Value = "isDeleted"
$ws.Range("H17").Value = "isDeleted"

$ws.Range("H24").Value = "categoryId"
$ws.Range("H25").Value = "authorId"
$ws.Range("H26").Value = "publisherId"
$ws.Range("H27").Value = "yearOfPublication"
$ws.Range("H31").Value = "isDeleted"

$ws.Range("H36").Value = "familyName"
$ws.Range("H37").Value = "givenName"
$ws.Range("H38").Value = "dateOfBirth"
$ws.Range("H42").Value = "isDeleted"

$ws.Range("H47").Value = "familyName"
$ws.Range("H48").Value = "givenName"
$ws.Range("H49").Value = "dateOfBirth"
$ws.Range("H53").Value = "citizenIdentification"
$ws.Range("H54").Value = "hashPassword"
$ws.Range("H57").Value = "isDeleted"

$ws.Range("H63").Value = "startDate"
$ws.Range("H64").Value = "endDate"

# Constraint note referencing the renamed "start_date" attribute
$ws.Range("J64").Value = "NOT NULL, > startDate"

$ws.Range("H66").Value = "discountPercent"
$ws.Range("H68").Value = "isDeleted"

$ws.Range("H73").Value = "staffId"
$ws.Range("H74").Value = "customerId"
$ws.Range("H75").Value = "promotionId"
$ws.Range("H76").Value = "createdTime"
$ws.Range("H77").Value = "totalAmount"
$ws.Range("H78").Value = "subTotalAmount"
$ws.Range("H79").Value = "promotionAmount"
$ws.Range("H82").Value = "isDeleted"

$ws.Range("H86").Value = "orderId"
$ws.Range("H87").Value = "bookId"
$ws.Range("H90").Value = "isDeleted"

# View-state tweaks captured in the saved workbook
$ws.Activate()
$ws.Range("J75").Select()
$excel.ActiveWindow.Zoom = 130
